$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 1.38  # G2: 1.44 -> 1.38
$ws.Cells.Item(2, 8).Value = 4.75  # H2: 4.33 -> 4.75
$ws.Cells.Item(2, 9).Value = 7  # I2: 6.5 -> 7
$ws.Cells.Item(2, 10).Value = 1.83  # J2: 1.95 -> 1.83
$ws.Cells.Item(2, 11).Value = 2.75  # K2: 2.6 -> 2.75
$ws.Cells.Item(2, 13).Value = 1.02  # M2: 1.03 -> 1.02
$ws.Cells.Item(2, 14).Value = 19  # N2: 17 -> 19
$ws.Cells.Item(2, 15).Value = 1.13  # O2: 1.14 -> 1.13
$ws.Cells.Item(2, 16).Value = 6  # P2: 5.5 -> 6
$ws.Cells.Item(2, 17).Value = 1.44  # Q2: 1.53 -> 1.44
$ws.Cells.Item(2, 18).Value = 2.7  # R2: 2.4 -> 2.7
$ws.Cells.Item(2, 19).Value = 1.22  # S2: 1.29 -> 1.22
$ws.Cells.Item(2, 20).Value = 4  # T2: 3.5 -> 4
$ws.Cells.Item(2, 21).Value = 1.62  # U2: 1.67 -> 1.62
$ws.Cells.Item(2, 22).Value = 2.2  # V2: 2.1 -> 2.2
$ws.Cells.Item(2, 23).Value = 11  # W2: 9 -> 11
$ws.Cells.Item(2, 25).Value = 9  # Y2: 8.5 -> 9
$ws.Cells.Item(2, 29).Value = 21  # AC2: 17 -> 21
$ws.Cells.Item(2, 30).Value = 10  # AD2: 8.5 -> 10
$ws.Cells.Item(2, 34).Value = 23  # AH2: 21 -> 23
$ws.Cells.Item(2, 40).Value = 3.75  # AN2: 3.6 -> 3.75
$ws.Cells.Item(2, 41).Value = 6.5  # AO2: 7 -> 6.5
$ws.Cells.Item(2, 43).Value = 17  # AQ2: 19 -> 17
$ws.Cells.Item(2, 44).Value = 34  # AR2: 41 -> 34
$ws.Cells.Item(2, 45).Value = 81  # AS2: 101 -> 81
$ws.Cells.Item(2, 46).Value = 4  # AT2: 3.5 -> 4
$ws.Cells.Item(2, 48).Value = 41  # AV2: 51 -> 41
$ws.Cells.Item(2, 49).Value = 8.5  # AW2: 8 -> 8.5
# Row 3
$ws.Cells.Item(3, 7).Value = 1.73  # G3: 1.7 -> 1.73
$ws.Cells.Item(3, 8).Value = 4  # H3: 3.9 -> 4
$ws.Cells.Item(3, 9).Value = 4.33  # I3: 4.5 -> 4.33
$ws.Cells.Item(3, 11).Value = 2.5  # K3: 2.4 -> 2.5
$ws.Cells.Item(3, 12).Value = 4.33  # L3: 4.5 -> 4.33
$ws.Cells.Item(3, 13).Value = 1.02  # M3: 1.03 -> 1.02
$ws.Cells.Item(3, 14).Value = 19  # N3: 17 -> 19
$ws.Cells.Item(3, 15).Value = 1.14  # O3: 1.17 -> 1.14
$ws.Cells.Item(3, 16).Value = 5.5  # P3: 5 -> 5.5
$ws.Cells.Item(3, 17).Value = 1.5  # Q3: 1.53 -> 1.5
$ws.Cells.Item(3, 18).Value = 2.5  # R3: 2.4 -> 2.5
$ws.Cells.Item(3, 19).Value = 1.25  # S3: 1.29 -> 1.25
$ws.Cells.Item(3, 20).Value = 3.75  # T3: 3.5 -> 3.75
$ws.Cells.Item(3, 21).Value = 1.5  # U3: 1.57 -> 1.5
$ws.Cells.Item(3, 22).Value = 2.5  # V3: 2.25 -> 2.5
$ws.Cells.Item(3, 23).Value = 11  # W3: 10 -> 11
$ws.Cells.Item(3, 24).Value = 11  # X3: 10 -> 11
$ws.Cells.Item(3, 28).Value = 19  # AB3: 21 -> 19
$ws.Cells.Item(3, 29).Value = 19  # AC3: 17 -> 19
$ws.Cells.Item(3, 31).Value = 12  # AE3: 13 -> 12
$ws.Cells.Item(3, 32).Value = 34  # AF3: 41 -> 34
$ws.Cells.Item(3, 33).Value = 101  # AG3: 126 -> 101
$ws.Cells.Item(3, 34).Value = 19  # AH3: 17 -> 19
$ws.Cells.Item(3, 39).Value = 29  # AM3: 34 -> 29
$ws.Cells.Item(3, 42).Value = 15  # AP3: 17 -> 15
$ws.Cells.Item(3, 45).Value = 81  # AS3: 101 -> 81
$ws.Cells.Item(3, 46).Value = 3.75  # AT3: 3.5 -> 3.75
$ws.Cells.Item(3, 47).Value = 7  # AU3: 7.5 -> 7
$ws.Cells.Item(3, 53).Value = 67  # BA3: 81 -> 67
$ws.Cells.Item(3, 55).Value = 301  # BC3: None -> 301
# Row 4
$ws.Cells.Item(4, 7).Value = 2.15  # G4: 2.2 -> 2.15
$ws.Cells.Item(4, 8).Value = 3.25  # H4: 3.2 -> 3.25
$ws.Cells.Item(4, 9).Value = 3.5  # I4: 3.4 -> 3.5
$ws.Cells.Item(4, 14).Value = 9.5  # N4: 10 -> 9.5
$ws.Cells.Item(4, 17).Value = 2  # Q4: 2.05 -> 2
$ws.Cells.Item(4, 18).Value = 1.8  # R4: 1.75 -> 1.8
$ws.Cells.Item(4, 19).Value = 1.4  # S4: 1.44 -> 1.4
$ws.Cells.Item(4, 20).Value = 2.75  # T4: 2.63 -> 2.75
$ws.Cells.Item(4, 21).Value = 1.73  # U4: 1.8 -> 1.73
$ws.Cells.Item(4, 22).Value = 2  # V4: 1.91 -> 2
$ws.Cells.Item(4, 25).Value = 9  # Y4: 9.5 -> 9
$ws.Cells.Item(4, 26).Value = 19  # Z4: 21 -> 19
$ws.Cells.Item(4, 27).Value = 17  # AA4: 19 -> 17
$ws.Cells.Item(4, 29).Value = 9.5  # AC4: 9 -> 9.5
$ws.Cells.Item(4, 33).Value = 201  # AG4: 251 -> 201
$ws.Cells.Item(4, 34).Value = 11  # AH4: 10 -> 11
$ws.Cells.Item(4, 36).Value = 13  # AJ4: 12 -> 13
$ws.Cells.Item(4, 37).Value = 41  # AK4: 34 -> 41
$ws.Cells.Item(4, 40).Value = 4  # AN4: 4.33 -> 4
$ws.Cells.Item(4, 42).Value = 21  # AP4: 23 -> 21
$ws.Cells.Item(4, 44).Value = 51  # AR4: 67 -> 51
$ws.Cells.Item(4, 46).Value = 2.75  # AT4: 2.63 -> 2.75
$ws.Cells.Item(4, 49).Value = 5.5  # AW4: 5 -> 5.5
$ws.Cells.Item(4, 52).Value = 67  # AZ4: 51 -> 67
# Row 5
$ws.Cells.Item(5, 7).Value = 2.1  # G5: 2.15 -> 2.1
$ws.Cells.Item(5, 9).Value = 3.5  # I5: 3.4 -> 3.5
$ws.Cells.Item(5, 12).Value = 3.75  # L5: 3.6 -> 3.75
$ws.Cells.Item(5, 24).Value = 11  # X5: 12 -> 11
$ws.Cells.Item(5, 26).Value = 19  # Z5: 21 -> 19
$ws.Cells.Item(5, 27).Value = 15  # AA5: 17 -> 15
$ws.Cells.Item(5, 36).Value = 13  # AJ5: 12 -> 13
$ws.Cells.Item(5, 37).Value = 41  # AK5: 34 -> 41
$ws.Cells.Item(5, 43).Value = 34  # AQ5: 41 -> 34
# Row 13
$ws.Cells.Item(13, 7).Value = 2.63  # G13: 2.6 -> 2.63
$ws.Cells.Item(13, 8).Value = 3.2  # H13: 3.3 -> 3.2
$ws.Cells.Item(13, 9).Value = 2.63  # I13: 2.6 -> 2.63
$ws.Cells.Item(13, 10).Value = 3.4  # J13: 3.25 -> 3.4
$ws.Cells.Item(13, 12).Value = 3.4  # L13: 3.25 -> 3.4
$ws.Cells.Item(13, 13).Value = 1.06  # M13: 1.05 -> 1.06
$ws.Cells.Item(13, 14).Value = 9.5  # N13: 11 -> 9.5
$ws.Cells.Item(13, 15).Value = 1.3  # O13: 1.29 -> 1.3
$ws.Cells.Item(13, 16).Value = 3.4  # P13: 3.5 -> 3.4
$ws.Cells.Item(13, 17).Value = 2.05  # Q13: 1.93 -> 2.05
$ws.Cells.Item(13, 18).Value = 1.75  # R13: 1.93 -> 1.75
$ws.Cells.Item(13, 19).Value = 1.44  # S13: 1.4 -> 1.44
$ws.Cells.Item(13, 20).Value = 2.63  # T13: 2.75 -> 2.63
$ws.Cells.Item(13, 21).Value = 1.8  # U13: 1.73 -> 1.8
$ws.Cells.Item(13, 22).Value = 1.91  # V13: 2 -> 1.91
$ws.Cells.Item(13, 23).Value = 8.5  # W13: 9 -> 8.5
$ws.Cells.Item(13, 29).Value = 9.5  # AC13: 10 -> 9.5
$ws.Cells.Item(13, 30).Value = 6  # AD13: 6.5 -> 6
$ws.Cells.Item(13, 31).Value = 15  # AE13: 13 -> 15
$ws.Cells.Item(13, 32).Value = 51  # AF13: 41 -> 51
$ws.Cells.Item(13, 33).Value = 251  # AG13: 201 -> 251
$ws.Cells.Item(13, 34).Value = 8.5  # AH13: 9 -> 8.5
$ws.Cells.Item(13, 42).Value = 26  # AP13: 23 -> 26
$ws.Cells.Item(13, 45).Value = 201  # AS13: 151 -> 201
$ws.Cells.Item(13, 46).Value = 2.63  # AT13: 2.75 -> 2.63
$ws.Cells.Item(13, 51).Value = 26  # AY13: 23 -> 26
$ws.Cells.Item(13, 54).Value = 201  # BB13: 151 -> 201
# Row 14
$ws.Cells.Item(14, 7).Value = 1.38  # G14: 1.3 -> 1.38
$ws.Cells.Item(14, 8).Value = 4.5  # H14: 5.25 -> 4.5
$ws.Cells.Item(14, 9).Value = 8  # I14: 9.5 -> 8
$ws.Cells.Item(14, 10).Value = 1.91  # J14: 1.8 -> 1.91
$ws.Cells.Item(14, 11).Value = 2.38  # K14: 2.5 -> 2.38
$ws.Cells.Item(14, 12).Value = 7.5  # L14: 8.5 -> 7.5
$ws.Cells.Item(14, 13).Value = 1.04  # M14: 1.03 -> 1.04
$ws.Cells.Item(14, 14).Value = 13  # N14: 15 -> 13
$ws.Cells.Item(14, 15).Value = 1.25  # O14: 1.2 -> 1.25
$ws.Cells.Item(14, 16).Value = 3.75  # P14: 4.33 -> 3.75
$ws.Cells.Item(14, 17).Value = 1.8  # Q14: 1.7 -> 1.8
$ws.Cells.Item(14, 18).Value = 2  # R14: 2.1 -> 2
$ws.Cells.Item(14, 19).Value = 1.36  # S14: 1.33 -> 1.36
$ws.Cells.Item(14, 20).Value = 3  # T14: 3.25 -> 3
$ws.Cells.Item(14, 21).Value = 2.2  # U14: 2.1 -> 2.2
$ws.Cells.Item(14, 22).Value = 1.62  # V14: 1.67 -> 1.62
$ws.Cells.Item(14, 23).Value = 6.5  # W14: 7 -> 6.5
$ws.Cells.Item(14, 26).Value = 8.5  # Z14: 8 -> 8.5
$ws.Cells.Item(14, 27).Value = 13  # AA14: 12 -> 13
$ws.Cells.Item(14, 28).Value = 34  # AB14: 29 -> 34
$ws.Cells.Item(14, 29).Value = 11  # AC14: 12 -> 11
$ws.Cells.Item(14, 30).Value = 9  # AD14: 10 -> 9
$ws.Cells.Item(14, 32).Value = 81  # AF14: 67 -> 81
$ws.Cells.Item(14, 34).Value = 17  # AH14: 23 -> 17
$ws.Cells.Item(14, 36).Value = 23  # AJ14: 29 -> 23
$ws.Cells.Item(14, 37).Value = 81  # AK14: 126 -> 81
$ws.Cells.Item(14, 38).Value = 51  # AL14: 67 -> 51
$ws.Cells.Item(14, 39).Value = 51  # AM14: 67 -> 51
$ws.Cells.Item(14, 41).Value = 6.5  # AO14: 6 -> 6.5
$ws.Cells.Item(14, 42).Value = 21  # AP14: 19 -> 21
$ws.Cells.Item(14, 43).Value = 19  # AQ14: 17 -> 19
$ws.Cells.Item(14, 46).Value = 3  # AT14: 3.25 -> 3
$ws.Cells.Item(14, 49).Value = 9  # AW14: 10 -> 9
# Row 16
$ws.Cells.Item(16, 7).Value = 2.38  # G16: 2.3 -> 2.38
$ws.Cells.Item(16, 9).Value = 3.3  # I16: 3.5 -> 3.3
$ws.Cells.Item(16, 10).Value = 3.25  # J16: 3.2 -> 3.25
$ws.Cells.Item(16, 11).Value = 1.83  # K16: 1.91 -> 1.83
$ws.Cells.Item(16, 13).Value = 1.11  # M16: 1.13 -> 1.11
$ws.Cells.Item(16, 14).Value = 6.5  # N16: 6 -> 6.5
$ws.Cells.Item(16, 21).Value = 2.2  # U16: 2.25 -> 2.2
$ws.Cells.Item(16, 22).Value = 1.62  # V16: 1.57 -> 1.62
$ws.Cells.Item(16, 23).Value = 6  # W16: 5.5 -> 6
$ws.Cells.Item(16, 24).Value = 10  # X16: 9.5 -> 10
$ws.Cells.Item(16, 25).Value = 11  # Y16: 10 -> 11
$ws.Cells.Item(16, 26).Value = 23  # Z16: 21 -> 23
$ws.Cells.Item(16, 31).Value = 19  # AE16: 21 -> 19
$ws.Cells.Item(16, 40).Value = 4.33  # AN16: 4 -> 4.33
$ws.Cells.Item(16, 42).Value = 34  # AP16: 29 -> 34
$ws.Cells.Item(16, 44).Value = 101  # AR16: 81 -> 101
$ws.Cells.Item(16, 51).Value = 34  # AY16: 41 -> 34
$ws.Cells.Item(16, 52).Value = 67  # AZ16: 81 -> 67
